$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E3"   = 16.508
    "E14"  = 16.95500000000001
    "E16"  = 16.505
    "E21"  = 16.559
    "E23"  = 16.55
    "E25"  = 17.493
    "E26"  = 16.42
    "E29"  = 16.801
    "E40"  = 16.541
    "E53"  = 16.613
    "E57"  = 16.675
    "E59"  = 16.507
    "E65"  = 16.774
    "E69"  = 17.563
    "E79"  = 17.06
    "E83"  = 16.569
    "E91"  = 17.612
    "E93"  = 17.258
    "E100" = 16.749
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
